# Weekly driver report update for 2025-04-28
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B14").Value = 338880
$ws.Range("B15").Value = 143869
$ws.Range("B19").Value = 68450
$ws.Range("B22").Value = 90508
$ws.Range("B24").Value = 52515
